$d = $word.ActiveDocument

# Locate the paragraph that ends with the "LOB1038..." requirement line.
# The deletion must start right after this paragraph (and its paragraph mark).
$startRange = $d.Content.Duplicate
[void]$startRange.Find.Execute(
    "LOB1038: Física Experimental I (Requisito fraco)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
[void]$startRange.Expand(4)   # wdParagraph - include the trailing paragraph mark
$startPos = $startRange.End

# Locate the end of the footer/copyright paragraph; deletion must stop right
# after this paragraph (and its paragraph mark), leaving the subsequent
# (empty / page-break) paragraphs untouched.
$endRange = $d.Content.Duplicate
[void]$endRange.Find.Execute(
    "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
[void]$endRange.Expand(4)   # wdParagraph - include the trailing paragraph mark
$endPos = $endRange.End

# Remove the empty paragraph + "Ver no Jupiter..." paragraph + the
# copyright/footer paragraph in one shot.
$deleteRange = $d.Range($startPos, $endPos)
$deleteRange.Delete()
